$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: update description text to mention interests ---
$ws.Range("B3").Value = "Upate user profile with interests"

# --- New row 4: S1_TC_T3 - get profile and verify interests ---
$ws.Range("A4").Value = "S1_TC_T3"
$ws.Range("B4").Value = "Get profile of the user and verify the name and interest /skill details."
$ws.Range("C4").Value = "1PPROFILE"
$ws.Range("D4").Value = "/users/user/(S1_TC_T1_truid)"
$ws.Range("E4").Value = "GET"
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").Style = "Normal"
$ws.Range("H4").Style = "Normal"
$ws.Range("I4").Value = "S1_TC_T1"
$ws.Range("J4").Value = "status=200||lastName=Yalamarthi||interest=computers"
$ws.Range("K4").Style = "Normal"
$ws.Range("L4").Value = "PASS"
$ws.Rows.Item(4).RowHeight = 45

# --- New row 5: S1_TC_T4 - update user image ---
$ws.Range("A5").Value = "S1_TC_T4"
$ws.Range("B5").Value = "Update user image"
$ws.Range("C5").Value = "1PPROFILE"
$ws.Range("D5").Value = "/users/user/(S1_TC_T1_truid)/image"
$ws.Range("E5").Value = "PUT"
$ws.Range("F5").Value = "Content-Type=application/json"
$ws.Range("G5").Style = "Normal"
$ws.Range("H5").Value = '{"imageContent":""}'
$ws.Range("I5").Value = "S1_TC_T1"
$ws.Range("J5").Value = "status=200"
$ws.Range("K5").Style = "Normal"
$ws.Range("L5").Value = "PASS"

# --- New row 6: S1_TC_T5 - get user image ---
$ws.Range("A6").Value = "S1_TC_T5"
$ws.Range("B6").Value = "Get user image"
$ws.Range("C6").Value = "1PPROFILE"
$ws.Range("D6").Value = "/users/user/(S1_TC_T1_truid)/image"
$ws.Range("E6").Value = "GET"
$ws.Range("F6").Style = "Normal"
$ws.Range("G6").Style = "Normal"
$ws.Range("H6").Style = "Normal"
$ws.Range("I6").Value = "S1_TC_T1"
$ws.Range("J6").Value = "status=200||imageContent="
$ws.Range("K6").Style = "Normal"
$ws.Range("L6").Value = "PASS"

# --- wrap-text styling on description / body cells ---
$ws.Range("B1").WrapText = $true
$ws.Range("B3").WrapText = $true
$ws.Range("B4").WrapText = $true
$ws.Range("B5").WrapText = $true
$ws.Range("B6").WrapText = $true
$ws.Range("H5").WrapText = $true

# --- column widths ---
$ws.Columns.Item(2).ColumnWidth = 29.592447916666668
$ws.Columns.Item(4).ColumnWidth = 48.451822916666664
$ws.Columns.Item(10).ColumnWidth = 51.166666666666664

# --- selection / view ---
$ws.Range("A6").Select()
